$wb = $excel.ActiveWorkbook

# Rename sheets: drop the "_updated" suffix
$wb.Worksheets.Item("Summary_updated").Name = "Summary"
$wb.Worksheets.Item("Attendance_updated").Name = "Attendance"
$wb.Worksheets.Item("Transfers_updated").Name = "Transfers"

# Work on the Transfers sheet: add a new "Status" column (G)
$ws = $wb.Worksheets.Item("Transfers")

# Copy header formatting from column F's header cell into the new G1 header cell
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("G1").Value = "Status"

# Match column width of column A (width 12) for the new column G
$ws.Range("G1").ColumnWidth = $ws.Range("A1").ColumnWidth

# Extend the autofilter range to include the new column
$ws.AutoFilterMode = $false
$ws.Range("A1:G1").AutoFilter() | Out-Null

# Update the hidden _FilterDatabase defined name to match the new filter range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Transfers!_FilterDatabase") {
        $n.RefersTo = "=Transfers!`$A`$1:`$G`$1"
    }
}
